# Applies the update described in the commit: "New evaluation of forces using
# substraction of data points instead of substraction of fit"
#
# Changes:
#  - Updates coefficient text strings in D5, D17, D19 on sheet "Frédéric"
#  - Updates computed statistics in E5:I5, E17:I17, E19:I19 on the same sheet
#  - Widens column D and scrolls/selects a different view on that sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Frédéric")

# --- Update coefficient text (shared strings) ---
$ws.Range("D5").Value = "[-0.006428638953316, 0.003532085496834, 0.167195264490172,6.660874387423239e-04,0.594410370645878,0.002368065163524]"
$ws.Range("D17").Value = "[-0.008111212221499   -0.024771353274546     -0.008297633291171    0.177246306723145}"
$ws.Range("D19").Value = "[-0.004010655576495     -0.000721901291023     -0.112383214091401     0.077819486465484]"

# --- Update numeric results for row 5 ---
$ws.Range("E5").Value = 0.306488942975441
$ws.Range("F5").Value = 0.01604872616039
$ws.Range("G5").Value = 1.51864619175753
$ws.Range("H5").Value = 0.079521096681094
$ws.Range("I5").Value = 19.0974

# --- Update numeric results for row 17 ---
$ws.Range("E17").Value = 0.323963710192727
$ws.Range("F17").Value = 0.045672975063812
$ws.Range("G17").Value = 1.02814195849866
$ws.Range("H17").Value = 0.144949266091047
$ws.Range("I17").Value = 7.093116

# --- Update numeric results for row 19 ---
$ws.Range("E19").Value = 0.91801336793276
$ws.Range("F19").Value = 0.088860928663791
$ws.Range("G19").Value = 2.47637296736596
$ws.Range("H19").Value = 0.23970544360762
$ws.Range("I19").Value = 10.3309

# --- Widen column D ---
# Excel's ColumnWidth (character units) differs from the OOXML <col width>
# by a constant offset (~5/6 of a character) in this engine's metrics, so
# back the literal 86 target out of that offset.
$ws.Columns.Item(4).ColumnWidth = 85.16666666666667

# --- Update the sheet view: scroll to D1 and select H23 ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("H23").Select()
